# Generate Report for Handback
# - Mark the 450c8980 row as "Handback transform failed" (was "Ready for handoff")
# - Record the handback/handoff file-name mismatch error detail for zh-cn and de-de
# - Widen the "Error Detail" column (P) to fit the new message

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"
$zhError = "Handback file name: b1nfmjig.tbj is different with handoff file name: 450c8980-c280-4c2d-a2ad-103ff8acd456.3c7d43010667e26a8e9aaccc8b108635855bbc0f.zh-cn."
$deError  = "Handback file name: b1nfmjig.tbj is different with handoff file name: 450c8980-c280-4c2d-a2ad-103ff8acd456.3c7d43010667e26a8e9aaccc8b108635855bbc0f.de-de."
# Column P (Error Detail) needs a stored/XML width of 40; this COM layer's
# ColumnWidth setter applies a constant +5/6 offset (quantized to 1/6) on
# write, so back that out here to land exactly on the target stored width.
$errorDetailColumnWidth = 39.166666666666664

# -- Overview sheet: status columns for the 450c8980 file (row 3) --
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# -- zh-cn sheet: Status (C) + Error Detail (P) for the 450c8980 row (row 3) --
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth

# -- de-de sheet: Status (C) + Error Detail (P) for the 450c8980 row (row 3) --
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth
